# Add the new "converted to €" summary table (rows 7-9) to the "riassunto"
# sheet, mirroring the existing k€-based table in rows 1-3 but expressed in
# plain € units (values multiplied by 1000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("riassunto")

# Row 7: same headers as row 1 (y1, y2, y3, y4)
$ws.Range("A7").Value = "y1"
$ws.Range("B7").Value = "y2"
$ws.Range("C7").Value = "y3"
$ws.Range("D7").Value = "y4"

# Row 8: values from row 2, rescaled from k€ to €  (x 10^3)
$ws.Range("A8").Formula = "=A2*10^3"
$ws.Range("B8").Formula = "=B2*10^3"
$ws.Range("C8:D8").Formula = "=C2*10^3"

# Row 9: unit labels, equivalent to row 3 but in € instead of k€
$ws.Range("A9").Value = "€"
$ws.Range("B9").Value = "€/MW"
$ws.Range("C9").Value = "€/km"
$ws.Range("D9").Value = "€/(MW*km)"

# Match the printer/page setup that Excel recorded for this sheet after the edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
